$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily mark the Price column as Text so COM does not
# silently coerce dotted price strings (e.g. "312.40") into
# floating point numbers when the Value is assigned below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.445.62"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.823.50"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "312.40"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.4242"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("D8").Value = "0.3617"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.07201"
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").Value = "0.8594"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").Value = "20.57"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "1.852.06"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "5.389"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").Value = "6.469"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").Value = "0.06925"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "80.26"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "0.000008880"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "27.481.64"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "5.127"
$ws.Range("E22").Value = "  +3.09%  "
$ws.Range("D23").Value = "10.93"
$ws.Range("E23").Value = "  +5.64%  "
$ws.Range("D24").Value = "2.079.52"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").Value = "1.985"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "154.93"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").Value = "5.140"
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("D29").Value = "114.17"
$ws.Range("E29").Value = "  -4.95%  "
$ws.Range("E30").Value = "  -4.04%  "
$ws.Range("D31").Value = "0.08831"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").Value = "0.7466"
$ws.Range("E32").Value = "  -2.39%  "
$ws.Range("D33").Value = "2.968"
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").Value = "4.529"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("D35").Value = "1.118"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "1.085"
$ws.Range("E37").Value = "  -1.40%  "
$ws.Range("D38").Value = "0.05276"
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "2.779"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").Value = "0.5060"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "0.1639"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").Value = "6.427"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D44").Value = "8.329"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").Value = "10.47"
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("D46").Value = "105.61"
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").Value = "0.4671"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "0.06441"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "1.611"
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("D51").Value = "63.58"
$ws.Range("E51").Value = "  -1.26%  "

# Restore the original (unformatted) look of the column now
# that the text values are safely stored.
$ws.Range("D2:D51").ClearFormats()
